$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" note text with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.39 = 50000.0 pesos`n✅ 50000.0 pesos = 12.38 = 975.57 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 80.7
$ws2.Range("O10").Value = 4035
$ws2.Range("N12").Value = 4040
$ws2.Range("O12").Value = 78.82599999999999
